# Weekly data refresh: a new week's worth of price data (2 rows: "Primera"
# and "Segunda" quality) is prepended into the data table at row 96,
# pushing every existing row below it down by two rows (the former rows
# 96-201 become rows 98-203; the dimension grows from A1:R201 to A1:R203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 96-97; everything from the old row 96
# onward (previously through row 201) shifts down to rows 98-203.
$ws.Rows("96:97").Insert()

# New row 96 - "Primera" quality entry for the latest date.
$ws.Range("A96").Value = 1
$ws.Range("B96").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C96").Value = "Arica y Parinacota"
$ws.Range("D96").Value = 44494
$ws.Range("E96").Value = 15
$ws.Range("F96").Value = 100112043
$ws.Range("G96").Value = "Pepino ensalada"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 120
$ws.Range("K96").Value = 5000
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = 5500
$ws.Range("N96").Value = "`$/caja 70 unidades"
$ws.Range("O96").Value = "Región de Arica y Parinacota"
$ws.Range("P96").Value = 79
$ws.Range("Q96").Value = 70
$ws.Range("R96").Value = "Hortaliza"

# New row 97 - "Segunda" quality entry for the same latest date.
$ws.Range("A97").Value = 1
$ws.Range("B97").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C97").Value = "Arica y Parinacota"
$ws.Range("D97").Value = 44494
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = 100112043
$ws.Range("G97").Value = "Pepino ensalada"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Segunda"
$ws.Range("J97").Value = 130
$ws.Range("K97").Value = 4000
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = 4500
$ws.Range("N97").Value = "`$/caja 100 unidades"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 45
$ws.Range("Q97").Value = 100
$ws.Range("R97").Value = "Hortaliza"
